# Apply the edit described by the diff:
# Insert a new data row at worksheet row 107 (pushing the existing rows
# 107..164 down to 108..165) and populate the new row with the values
# for the added daily price record. All other rows keep their original
# content, simply shifted down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 107; this shifts rows
# 107-164 down to 108-165 and carries formatting from the row above,
# which matches the date-formatted style (s="2") used in column D.
$ws.Rows.Item(107).Insert()

# Populate the newly inserted row 107 with the new record's values.
$ws.Cells.Item(107, 1).Value  = 11
$ws.Cells.Item(107, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(107, 3).Value  = "Bíobío"
$ws.Cells.Item(107, 4).Value  = 45001
$ws.Cells.Item(107, 5).Value  = 8
$ws.Cells.Item(107, 6).Value  = 100112024
$ws.Cells.Item(107, 7).Value  = "Choclo"
$ws.Cells.Item(107, 8).Value  = "Choclero"
$ws.Cells.Item(107, 9).Value  = "Primera"
$ws.Cells.Item(107, 10).Value = 15000
$ws.Cells.Item(107, 11).Value = 400
$ws.Cells.Item(107, 12).Value = 450
$ws.Cells.Item(107, 13).Value = 417
$ws.Cells.Item(107, 14).Value = "$/unidad"
$ws.Cells.Item(107, 15).Value = "Región Metropolitana"
$ws.Cells.Item(107, 16).Value = 417
$ws.Cells.Item(107, 17).Value = 1
$ws.Cells.Item(107, 18).Value = "Hortaliza"
